$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Mit csináltam ma" (What I did today) header in C2,
# matching the existing header row pattern (Nevek / Mit csináltam ma / Mit fogok csinálni legközelebb / Van-e akadály?)
$ws.Range("C2").Value = "Mit csináltam ma"

# Move the active selection to C11, matching the saved cursor position
$ws.Range("C11").Select() | Out-Null
